$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.800.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +9.30%  "

$ws.Range("D3").Value = "'3.466.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.01%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'413.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.81%  "

$ws.Range("D6").Value = "'123.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.64%  "

$ws.Range("D7").Value = "'3.455.30"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.75%  "

$ws.Range("D8").Value = "'0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.89%  "

$ws.Range("E10").Value = "  +9.92%  "

$ws.Range("E11").Value = "  +33.81%  "

$ws.Range("D12").Value = "'41.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.80%  "

$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("D14").Value = "'4.011.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.03%  "

$ws.Range("D15").Value = "'8.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.46%  "

$ws.Range("D16").Value = "'19.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.03%  "

$ws.Range("D17").Value = "'3.459.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.80%  "

$ws.Range("D18").Value = "'62.758.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.59%  "

$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").Value = "'10.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.53%  "

$ws.Range("D21").Value = "'0.0000138"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +28.31%  "

$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Value = "'315.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.15%  "

$ws.Range("D24").Value = "'81.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.80%  "

$ws.Range("D25").Value = "'12.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("D26").Value = "'3.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.99%  "

$ws.Range("D27").Value = "'30.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.79%  "

$ws.Range("D28").Value = "'7.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.65%  "

$ws.Range("D29").Value = "'7.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.82%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.175"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.22%  "

$ws.Range("B31").Value = "LEO"
$ws.Range("C31").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D31").Value = "'4.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.01%  "

$ws.Range("D32").Value = "'0.117"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.33%  "

$ws.Range("D33").Value = "'2.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +23.10%  "

$ws.Range("D34").Value = "'11.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.76%  "

$ws.Range("D35").Value = "'42.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.10%  "

$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("D37").Value = "'0.0494"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.08%  "

$ws.Range("D38").Value = "'52.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.38%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'3.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.40%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.23%  "

$ws.Range("E41").Value = "  -2.13%  "

$ws.Range("E42").Value = "  +6.15%  "

$ws.Range("E43").Value = "  +3.16%  "

$ws.Range("D44").Value = "'135.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.56%  "

$ws.Range("D45").Value = "'0.283"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").Value = "'16.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("E47").Value = "  -0.43%  "

$ws.Range("E48").Value = "  +2.06%  "

$ws.Range("D49").Value = "'21.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.78%  "

$ws.Range("D50").Value = "'2.206.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.25%  "

$ws.Range("E51").Value = "  +0.38%  "
